$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing Latvian translations (currently in column A, rows 1-12)
$latvian = @(
    "Iegarena",
    "Tendrils",
    "lodveida kopas",
    "plūdmaiņu straumes",
    "Tuvā infrasarkanā kamera",
    "tumšā matērija",
    "liela mēroga struktūra",
    "Gadījumos",
    "viela",
    "deflācija",
    "Oļi",
    "hematīts"
)

# New English words to place in column A, rows 2-13
$english = @(
    "elongated",
    "tendrils",
    "globular clusters",
    "tidal streams",
    "near-infrared camera",
    "dark matter",
    "large-scale structure",
    "instances",
    "substance",
    "deflating",
    "pebbles",
    "hematite"
)

# Header row
$ws.Cells.Item(1, 1).Value = "English Word"
$ws.Cells.Item(1, 2).Value = "Latvian Translation"

# Move the Latvian words into column B, and set English words into column A
for ($i = 0; $i -lt $latvian.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $english[$i]
    $ws.Cells.Item($row, 2).Value = $latvian[$i]
}
